# "ideas table and meeting notes"
#
# Highlights two cells (first column) of the dataset-summary table in
# yellow:
#   1. "Wikipedia Abusive Language Data Set" row - the run AND the
#      paragraph mark get the yellow highlight.
#   2. "Civil Comments Toxicity Kaggle (CCTK)" row - only the four runs
#      that make up the cell text get the yellow highlight (no
#      paragraph-mark formatting is introduced here).

$wdYellowHighlight = 7

$d = $word.ActiveDocument

# The document has a single table (the dataset summary); locate it
# defensively in case other tables are ever introduced upstream.
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Columns.Count -ge 1) {
        $table = $candidate
        break
    }
}
if ($table -eq $null) {
    $table = $d.Tables.Item(1)
}

$wikipediaRow = 0
$civilCommentsRow = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $cellText = $table.Cell($r, 1).Range.Text
    if ($wikipediaRow -eq 0 -and $cellText -like "*Wikipedia Abusive Language Data Set*") {
        $wikipediaRow = $r
    }
    if ($civilCommentsRow -eq 0 -and $cellText -like "*Civil Comments*Toxicity*Kaggle*CCTK*") {
        $civilCommentsRow = $r
    }
}

# --- 1. "Wikipedia Abusive Language Data Set" --------------------------
# Apply highlight to the whole paragraph (Paragraph.Range), which covers
# both the run's rPr and the trailing paragraph mark's pPr/rPr, matching
# how Word applies formatting when the paragraph (incl. pilcrow) is
# selected.
if ($wikipediaRow -gt 0) {
    $wikiCell = $table.Cell($wikipediaRow, 1)
    $wikiPara = $wikiCell.Range.Paragraphs.Item(1)
    $wikiPara.Range.Font.HighlightColorIndex = $wdYellowHighlight
}

# --- 2. "Civil Comments Toxicity Kaggle (CCTK)" -------------------------
# Apply highlight only to the cell's text runs (not the paragraph mark),
# using a plain text Range so no pPr/rPr gets introduced.
if ($civilCommentsRow -gt 0) {
    $ccCell = $table.Cell($civilCommentsRow, 1)
    $ccStart = $ccCell.Range.Start
    $ccEnd = $ccCell.Range.End
    $ccRange = $d.Range($ccStart, $ccEnd)
    $ccRange.Font.HighlightColorIndex = $wdYellowHighlight
}
